$wb = $excel.ActiveWorkbook

# 1. Duplicate the "固件寄存器格式" sheet (index 3) right after itself.
$wsSrc = $wb.Worksheets.Item(3)
$wsSrc.Copy([System.Reflection.Missing]::Value, $wsSrc)
$ws = $wb.Worksheets.Item(4)
$ws.Name = "固件寄存器格式 - 修改"

# 2. Insert a brand-new blank column E so the register table can show both
#    the I-DAC and Q-DAC data fields side by side.
$ws.Columns.Item(5).Insert()

# 3. Row 1: "LSB" moves from D1 into the new E1 header cell.
$ws.Range("E1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = ""

# 4. Rows 2-4 (reg0 / DAC_CTL): merge D:E so the control-bit descriptions
#    keep spanning the full width, and center them like before.
$ws.Range("D2:E2").Merge()
$ws.Range("D2:E2").HorizontalAlignment = -4108
$ws.Range("D3:E3").Merge()
$ws.Range("D3:E3").HorizontalAlignment = -4108
$ws.Range("D4:E4").Merge()
$ws.Range("D4:E4").HorizontalAlignment = -4108

# 5. Tidy up the FRMT_CTL description text (drop the embedded line break).
$ws.Range("C4").Value = "格式控制位，0表示无符号二进制数，1表示二进制补码，默认为0"

# 6. Row 5-7 (reg1 / DAC_DATA): split the old merged B:C "DAC_DATABITS"
#    column into separate Q-DAC (B/C) and I-DAC (D/E) columns.
$ws.Range("B5:C5").UnMerge()
$ws.Range("B6:C6").UnMerge()
$ws.Range("B7:C7").UnMerge()

# The old b31-b10 / N/A / DAC_DATABITS / DAC数据位 fields become the Q-DAC
# columns (B/C), while the existing b9-b0 field shifts right to make room
# for the new I-DAC bit-range column (D).
$ws.Range("E5").Value = $ws.Range("D5").Value()
$ws.Range("D5").Value = "b15-b10"
$ws.Range("B5").Value = "b31-b26"
$ws.Range("C5").Value = "b25-b16"

$ws.Range("E6").Value = "DAC_I_DATA"
$ws.Range("D6").Value = "N/A"
$ws.Range("C6").Value = "DAC_Q_DATA"

$ws.Range("E7").Value = "I DAC数据位"
$ws.Range("D7").Value = "备用"
$ws.Range("C7").Value = "Q DAC数据位"

# Reset the alignment on the now-unmerged B/C cells back to "general".
$ws.Range("B5").HorizontalAlignment = 1
$ws.Range("C5").HorizontalAlignment = 1
$ws.Range("B6").HorizontalAlignment = 1
$ws.Range("C6").HorizontalAlignment = 1
$ws.Range("B7").HorizontalAlignment = 1
$ws.Range("C7").HorizontalAlignment = 1

# 7. Column widths for the new layout.
$ws.Columns.Item(3).ColumnWidth = 31.7109375
$ws.Columns.Item(4).ColumnWidth = 8.5703125
$ws.Columns.Item(5).ColumnWidth = 31.28515625

# 8. Selections: land on the new sheet with the whole table selected, and
#    leave the original sheet's selection parked away from the table.
$ws.Range("A1:E7").Select()
$wsSrc.Range("B10").Select()
$ws.Activate()
